$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns.Item(7).ColumnWidth = 10.99
